$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("_id", "isActive", "employeeDetails", "leavedetails", "employeeOfficeDetails", "month", "year", "leave_type", "leaveTypeName", "leaveStatus")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$employeeDetailsJson = '{"_id":2,"userName":"1010002","fullName":"Test Group HR Head 2"}'
$leaveDetailsJson = '{"_id":1,"updatedAt":"2018-08-17T05:26:15.665Z","createdAt":"2018-06-24T08:38:43.999Z","leave_type":1,"createdBy":8,"updatedBy":2,"attachment":"externalDocument/a42d4cd4.api list-2.pdf","status":"Pending Withdrawal","reason":"Rejected","days":2,"toDate":"2018-08-14T07:53:16.000Z","fromDate":"2018-08-12T08:53:16.000Z","createdByName":{"_id":8,"fullName":"Akshay k"},"updatedByName":{"_id":2,"fullName":"Test Group HR Head 2"},"supervisorDetails":[{"_id":5,"fullName":"Reviewer 5"}]}'
$employeeOfficeDetailsJson = '{"_id":2,"emp_id":2,"departments":{"_id":1,"departmentName":"International Business"},"divisions":{"_id":1,"divisionName":"Business"}}'

# Data row
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = $true
$ws.Cells.Item(2, 3).Value = $employeeDetailsJson
$ws.Cells.Item(2, 4).Value = $leaveDetailsJson
$ws.Cells.Item(2, 5).Value = $employeeOfficeDetailsJson
$ws.Cells.Item(2, 6).Value = 8
$ws.Cells.Item(2, 7).Value = 2018
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = "Annual Leave"
$ws.Cells.Item(2, 10).Value = "Pending Withdrawal"

Write-Host "done"
